$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# "Generate Report for Handoff" - a new handoff batch was generated for
# b.md: it now has a freshly generated xliff handoff file/date, its
# status flips from "Handed back: in sync with en-US" to "Ready for
# handoff", the Content Duplicate flag flips to False and an Error
# Detail note is recorded because the handback file on file is stale.
# ----------------------------------------------------------------------

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7788ff499055aaffebce5157132284b3b540754c/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/59e81fc6ddc20b7277cfeb80d9d664a26989a7b0/e2e/b.md."

# --- Overview sheet: row for b.md ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-23 12:37:15"

# --- zh-cn sheet: row for b.md ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Ready for handoff"
# F3 needs the literal text "False" (not a Boolean). F2 in this same
# column already stores "False" as text, so copy/paste-values from it
# instead of re-typing the literal (which Excel auto-converts to a
# Boolean value when assigned directly via .Value).
$wsZh.Range("F2").Copy()
$wsZh.Range("F3").PasteSpecial(-4163) # xlPasteValues
$wsZh.Application.CutCopyMode = $false
$wsZh.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-23 12:37:09"
$wsZh.Range("P3").Value = $errorDetail
$wsZh.Columns.Item(16).ColumnWidth = 40

# --- de-de sheet: row for b.md ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("F2").Copy()
$wsDe.Range("F3").PasteSpecial(-4163) # xlPasteValues
$wsDe.Application.CutCopyMode = $false
$wsDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-23 12:37:15"
$wsDe.Range("P3").Value = $errorDetail
$wsDe.Columns.Item(16).ColumnWidth = 40
